$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output "NOT FOUND: $old"
    }
}

# ------------------------------------------------------------------
# The edit rotates several blocks of text between paragraphs (the
# paragraph styles/formatting stay put; only the w:t content moves
# around). To perform this safely with simple Find/Replace we first
# swap every source string out for a unique placeholder token, then
# swap every placeholder token in for its final destination text.
#
# Rotation (content letter -> destination slot):
#   A (Objetivos/PT body)        -> Docente(s) value
#   B (Objetivos/EN body)        -> Programa resumido / EN body
#   C (Docente value)            -> Bibliografia value
#   D (Programa resumido/PT)     -> Objetivos / PT body
#   E (Programa resumido/EN)     -> Objetivos / EN body
#   F (Programa/PT body)         -> Programa resumido / PT body
#   G (Metodo value)             -> Programa / PT body
#   H (Criterio value)           -> Metodo value
#   I (Norma de recuperacao val) -> Criterio value
#   J (Bibliografia value)       -> Norma de recuperacao value
# ------------------------------------------------------------------

$br = [char]11
$biblio = "1. S. Zachs, INTRODUCTION TO RELIABILITY ANALYSIS: PROBABILITY MODELS AND STATISTICAL METHODS, Springer Verlag, New York, 1992" + $br + `
          "2. I.B. Gertsbakh, STATISTICAL RELIABILITY THEORY, Marcel Dekker, New York, 1989." + $br + `
          "3. J. Knezevic, RELIABILITY, MAINTAINABILITY, AND SUPPORTABILITY: A PROBABILITY APPROACH, McGraw-Hill, 1993." + $br + `
          "4. R.S. Dhillon, C. Singh, ENGINEERING RELIABILITY. NEW TECHNIQUES AND APPLICATIONS, Wiley Interscience, 1981. " + $br + `
          "5. HARRY, M. , LINSENMANND.R., The Six Sigma Fieldbook, Doubleday, New York, 2006"

$A = "Explicar os conceitos, métodos e resolver problemas que ilustrem aplicações sem recorrer a desenvolvimento teóricos da Teoria da Confiabilidade. Pretende-se uma formação geral com o uso de modelos probabilísticos e estatísticos, e com aplicações na área de engenharia. Uso de aplicativos computacionais para análise de conjunto de dados."
$B = "Explain the concepts, methods and solve problems that illustrate applications without resorting to theoretical development of Reliability Theory. It is intended a general formation with the use of probabilistic and statistical models, and with applications in the area of engineering. Use of computational applications for data set analysis."
$C = "3295113 - José Eduardo Holler Branco"
$D = "1. Confiabilidade e disponibilidade de sistemas. 2. Famílias de distribuições. 3. Sistemas reparáveis. 4. Análise gráfica de dados. 5. Estimação de características de um sistema. 6. FMEA, 7. Aplicações na Gestão da Manutenção, 8. Manutenção Produtiva Total, 9. Design for Six Sigma. 10. RCM (Reliability Centered Maintenance)"
$E = "1. Reliability and availability of systems. 2. Families of distributions. 3. Repairable systems. 4. Graphical analysis of data. 5. Estimation of characteristics of a system. 6. FMEA, 7. Maintenance Management Applications, 8. Total Productive Maintenance, 9. Design for Six Sigma. 10. RCM (Reliability Centered Maintenance )"
$F = "1. Confiabilidade e disponibilidade de sistemas, decomposição por cortes e caminhos, árvores de eventos. 2. Famílias de distribuições úteis em Teoria a Confiabilidade. 3. Sistemas reparáveis, manutenção, aproximações assintóticas. 4. Análise gráfica de dados. 5. Estimação do tempo de vida e das características de um sistema. 6. FMEA, 7. Aplicações na Gestão da Manutenção, 8. Manutenção Produtiva Total, 9. Design for Six Sigma DFSS 10. RCM (Reliability Centered Maintenance)"
$G = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$H = "NF≥ 5,0."
$I = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
$J = $biblio

# --- Step 1: move everything currently in the doc out to unique placeholders ---

Replace-Text $A "@@TOKEN_A@@"
Replace-Text $B "@@TOKEN_B@@"
Replace-Text $C "@@TOKEN_C@@"
Replace-Text $D "@@TOKEN_D@@"
Replace-Text $E "@@TOKEN_E@@"
Replace-Text $F "@@TOKEN_F@@"
Replace-Text $G "@@TOKEN_G@@"
Replace-Text $H "@@TOKEN_H@@"
Replace-Text $I "@@TOKEN_I@@"
Replace-Text $J "@@TOKEN_J@@"

# --- Step 2: move placeholders into their final destinations ----------
# Each placeholder currently sits in the paragraph slot that used to hold
# that content; replace it with whatever content should *now* live there.

Replace-Text "@@TOKEN_A@@" $D   # slot "Objetivos"/PT now holds D
Replace-Text "@@TOKEN_B@@" $E   # slot "Objetivos"/EN now holds E
Replace-Text "@@TOKEN_C@@" $A   # slot "Docente(s)" value now holds A
Replace-Text "@@TOKEN_D@@" $F   # slot "Programa resumido"/PT now holds F
Replace-Text "@@TOKEN_E@@" $B   # slot "Programa resumido"/EN now holds B
Replace-Text "@@TOKEN_F@@" $G   # slot "Programa"/PT now holds G
Replace-Text "@@TOKEN_G@@" $H   # slot "Método" value now holds H
Replace-Text "@@TOKEN_H@@" $I   # slot "Critério" value now holds I
Replace-Text "@@TOKEN_I@@" $J   # slot "Norma de recuperação" value now holds J
Replace-Text "@@TOKEN_J@@" $C   # slot "Bibliografia" value now holds C

Write-Output "done"
